$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "WithTable" (contains ListObject "MyTable", range A1:E5)
# Insert a new "Boolean" column as the 4th column (between "DateTime" and
# "ARCtrl Column"), shifting the two "ARCtrl Column" columns one place right.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("WithTable")
$lo1 = $ws1.ListObjects.Item(1)

# Grow the table by one column (new column appears at the end, range F1:F5).
$lo1.Resize($ws1.Range("A1:F5"))

# Shift the old D:E columns (ARCtrl Column / ARCtrl Column<space>) into E:F.
$ws1.Range("E1:F5").Value = $ws1.Range("D1:E5").Value()

# Clear out D and populate it with the new Boolean column.
$ws1.Range("D1:D5").ClearContents()
$ws1.Range("D1").Value = "Boolean"
$ws1.Range("D2").Value = $true
$ws1.Range("D3").Value = $false
$ws1.Range("D4").Value = $true
$ws1.Range("D5").Value = $false

# Give the new boolean cells the same cell format used by the DateTime column.
$ws1.Range("C2:C5").Copy()
$ws1.Range("D2:D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Tableless" - same logical layout as sheet 1, but no ListObject.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Tableless")

$ws2.Range("E1:F5").Value = $ws2.Range("D1:E5").Value()

$ws2.Range("D1:D5").ClearContents()
$ws2.Range("D1").Value = "Boolean"
$ws2.Range("D2").Value = $true
$ws2.Range("D3").Value = $false
$ws2.Range("D4").Value = $true
$ws2.Range("D5").Value = $false

$ws2.Range("C2:C5").Copy()
$ws2.Range("D2:D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give column D an explicit width like column C's.
$ws2.Columns.Item(4).ColumnWidth = 9.33

# ---------------------------------------------------------------------------
# Sheet 3: "WithTable_Duplicate" (contains ListObject "MyOtherTable",
# range B4:F8). Same insert, offset by one column (table starts at B).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("WithTable_Duplicate")
$lo3 = $ws3.ListObjects.Item(1)

$lo3.Resize($ws3.Range("B4:G8"))

$ws3.Range("F4:G8").Value = $ws3.Range("E4:F8").Value()

$ws3.Range("E4:E8").ClearContents()
$ws3.Range("E4").Value = "Boolean"
$ws3.Range("E5").Value = $true
$ws3.Range("E6").Value = $false
$ws3.Range("E7").Value = $true
$ws3.Range("E8").Value = $false

$ws3.Range("D5:D8").Copy()
$ws3.Range("E5:E8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selections / active sheet. "WithTable" ends up the active (first) tab,
# "Tableless" loses its tab selection, and each sheet keeps its own
# last-used selection rectangle.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F1").Select()

$ws3.Activate()
$ws3.Range("F7").Select()

$ws1.Activate()
$ws1.Range("D9").Select()
